$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the date/volume/price figures between row 2 and row 4
# (re-ordering the weekly entries into chronological order).

# Row 2 -> new values (previously row 4's values)
$ws.Range("D2").Value = 44284
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 500

# Row 4 -> new values (previously row 2's values)
$ws.Range("D4").Value = 44277
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550
